$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.245.56'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '2.646.93'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''597.20'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').Value = '''156.71'
$ws.Range('E6').Value = '  +1.29%  '
$ws.Range('E9').Value = '  +2.88%  '
$ws.Range('E10').Value = '  -1.19%  '
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('E14').Value = '  +1.75%  '
$ws.Range('D15').Value = '3.128.53'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').Value = '68.204.60'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').Value = '2.656.52'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('D19').Value = '''363.46'
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('E20').Value = '  -1.15%  '
$ws.Range('E21').Value = '  +3.40%  '
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('E23').Value = '  -2.48%  '
$ws.Range('D24').Value = '''75.02'
$ws.Range('E24').Value = '  +2.13%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').Value = '''9.77'
$ws.Range('E26').Value = '  -2.50%  '
$ws.Range('D27').Value = '2.779.37'
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('E28').Value = '  -0.58%  '
$ws.Range('D30').Value = '''559.72'
$ws.Range('E30').Value = '  -2.68%  '
$ws.Range('D31').Value = '''8.05'
$ws.Range('E31').Value = '  +0.61%  '
$ws.Range('E32').Value = '  -1.06%  '
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('D37').Value = '''161.69'
$ws.Range('E37').Value = '  +1.49%  '
$ws.Range('D38').Value = '''19.67'
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('E39').Value = '  +1.13%  '
$ws.Range('E40').Value = '  -2.67%  '
$ws.Range('D41').Value = '''5.32'
$ws.Range('E41').Value = '  -1.16%  '
$ws.Range('E42').Value = '  +3.54%  '
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('E44').Value = '  -1.43%  '
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').Value = '''158.86'
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').Value = '''22.06'
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('E51').Value = '  -0.24%  '
